$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newTickers = @("IMX-USD", "MNT-USD", "TAO-USD")

$startRow = 270
for ($i = 0; $i -lt $newTickers.Length; $i++) {
    $ws.Cells.Item($startRow + $i, 1).Value = $newTickers[$i]
}
